# Auto-generated Excel COM-interop script
# Applies numeric-value updates to the Chocobo_Profits workbook's per-sheet leve tables
# (columns H/I/J/K/L/M/N = currentAveragePrice*, LevePrice*, LeveProfit*).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5925
$ws.Range("J40").Value = 6342.857
$ws.Range("L40").Value = 6342.857
$ws.Range("N40").Value = -6692.857
# Row 124
$ws.Range("H124").Value = 42000
$ws.Range("J124").Value = 42000
$ws.Range("L124").Value = 42000
$ws.Range("N124").Value = -51820
# Row 141
$ws.Range("H141").Value = 5349.1787
$ws.Range("I141").Value = 5540.7085
$ws.Range("J141").Value = 4200
$ws.Range("K141").Value = 16622.1255
$ws.Range("L141").Value = 12600
$ws.Range("M141").Value = -11442.1255
$ws.Range("N141").Value = -22960

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2667.0908
$ws.Range("I74").Value = 2029.3529
$ws.Range("J74").Value = 4835.4
$ws.Range("K74").Value = 2029.3529
$ws.Range("L74").Value = 4835.4
$ws.Range("M74").Value = -1155.3529
$ws.Range("N74").Value = -6583.4
# Row 77
$ws.Range("H77").Value = 2667.0908
$ws.Range("I77").Value = 2029.3529
$ws.Range("J77").Value = 4835.4
$ws.Range("K77").Value = 10146.7645
$ws.Range("L77").Value = 24177
$ws.Range("M77").Value = -5778.764500000001
$ws.Range("N77").Value = -32913
# Row 137
$ws.Range("H137").Value = 40881.125
$ws.Range("J137").Value = 40881.125
$ws.Range("L137").Value = 40881.125
$ws.Range("N137").Value = -51081.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2472.1177
$ws.Range("J99").Value = 3512
$ws.Range("L99").Value = 3512
$ws.Range("N99").Value = -6508
# Row 137
$ws.Range("H137").Value = 39576
$ws.Range("J137").Value = 40640
$ws.Range("L137").Value = 40640
$ws.Range("N137").Value = -50840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 22499
$ws.Range("I12").Value = 19998
$ws.Range("K12").Value = 19998
$ws.Range("M12").Value = -19828
# Row 31
$ws.Range("H31").Value = 16132212
$ws.Range("I31").Value = 1899.0416
$ws.Range("K31").Value = 1899.0416
$ws.Range("M31").Value = -1604.0416
# Row 34
$ws.Range("H34").Value = 16132212
$ws.Range("I34").Value = 1899.0416
$ws.Range("K34").Value = 1899.0416
$ws.Range("M34").Value = -1697.0416
# Row 132
$ws.Range("H132").Value = 3125.8647
$ws.Range("I132").Value = 1786.9474
$ws.Range("J132").Value = 4539.1665
$ws.Range("K132").Value = 5360.8422
$ws.Range("L132").Value = 13617.4995
$ws.Range("M132").Value = -2830.8422
$ws.Range("N132").Value = -18677.4995
# Row 134
$ws.Range("H134").Value = 6491.04
$ws.Range("I134").Value = 8863.538
$ws.Range("J134").Value = 3920.8333
$ws.Range("K134").Value = 26590.614
$ws.Range("L134").Value = 11762.4999
$ws.Range("M134").Value = -24055.614
$ws.Range("N134").Value = -16832.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3483.75
$ws.Range("I3").Value = 2306.1538
$ws.Range("J3").Value = 8586.666999999999
$ws.Range("K3").Value = 6918.4614
$ws.Range("L3").Value = 25760.001
$ws.Range("M3").Value = -6806.4614
$ws.Range("N3").Value = -25984.001
# Row 87
$ws.Range("H87").Value = 4000
$ws.Range("I87").Value = 4000
$ws.Range("K87").Value = 12000
$ws.Range("M87").Value = -10752
# Row 90
$ws.Range("H90").Value = 4000
$ws.Range("I90").Value = 4000
$ws.Range("K90").Value = 36000
$ws.Range("M90").Value = -29760
# Row 136
$ws.Range("H136").Value = 3637.7727
$ws.Range("I136").Value = 3126.2307
$ws.Range("J136").Value = 4376.6665
$ws.Range("K136").Value = 9378.6921
$ws.Range("L136").Value = 13129.9995
$ws.Range("M136").Value = -4278.6921
$ws.Range("N136").Value = -23329.9995
# Row 137
$ws.Range("H137").Value = 3071.2222
$ws.Range("I137").Value = 765
$ws.Range("J137").Value = 3359.5
$ws.Range("K137").Value = 2295
$ws.Range("L137").Value = 10078.5
$ws.Range("M137").Value = 2805
$ws.Range("N137").Value = -20278.5
# Row 140
$ws.Range("H140").Value = 22288.54
$ws.Range("I140").Value = 38683.715
$ws.Range("J140").Value = 3160.8333
$ws.Range("K140").Value = 116051.145
$ws.Range("L140").Value = 9482.499899999999
$ws.Range("M140").Value = -110871.145
$ws.Range("N140").Value = -19842.4999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576
# Row 81
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
# Row 84
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
# Row 88
$ws.Range("H88").Value = 32487.5
$ws.Range("J88").Value = 32487.5
$ws.Range("L88").Value = 32487.5
$ws.Range("N88").Value = -33389.5
# Row 91
$ws.Range("H91").Value = 32487.5
$ws.Range("J91").Value = 32487.5
$ws.Range("L91").Value = 32487.5
$ws.Range("N91").Value = -35607.5
# Row 137
$ws.Range("H137").Value = 58298.707
$ws.Range("J137").Value = 58942.375
$ws.Range("L137").Value = 58942.375
$ws.Range("N137").Value = -69142.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 762.0682
$ws.Range("I68").Value = 710.02325
$ws.Range("K68").Value = 710.02325
$ws.Range("M68").Value = 38.97675000000004
# Row 71
$ws.Range("H71").Value = 762.0682
$ws.Range("I71").Value = 710.02325
$ws.Range("K71").Value = 3550.11625
$ws.Range("M71").Value = 193.88375
# Row 80
$ws.Range("H80").Value = 43025.6
$ws.Range("J80").Value = 43025.6
$ws.Range("L80").Value = 43025.6
$ws.Range("N80").Value = -45271.6
# Row 82
$ws.Range("H82").Value = 1610.826
$ws.Range("I82").Value = 771.8125
$ws.Range("J82").Value = 3528.5715
$ws.Range("K82").Value = 771.8125
$ws.Range("L82").Value = 3528.5715
$ws.Range("M82").Value = -410.8125
$ws.Range("N82").Value = -4250.5715
# Row 83
$ws.Range("H83").Value = 43025.6
$ws.Range("J83").Value = 43025.6
$ws.Range("L83").Value = 129076.8
$ws.Range("N83").Value = -140308.8
# Row 85
$ws.Range("H85").Value = 1610.826
$ws.Range("I85").Value = 771.8125
$ws.Range("J85").Value = 3528.5715
$ws.Range("K85").Value = 771.8125
$ws.Range("L85").Value = 3528.5715
$ws.Range("M85").Value = 476.1875
$ws.Range("N85").Value = -6024.5715
# Row 86
$ws.Range("H86").Value = 34130
$ws.Range("J86").Value = 34130
$ws.Range("L86").Value = 34130
$ws.Range("N86").Value = -36502
# Row 89
$ws.Range("H89").Value = 34130
$ws.Range("J89").Value = 34130
$ws.Range("L89").Value = 102390
$ws.Range("N89").Value = -114246
# Row 136
$ws.Range("H136").Value = 2445.641
$ws.Range("I136").Value = 1363.0303
$ws.Range("K136").Value = 4089.0909
$ws.Range("M136").Value = -1539.0909

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10205094
$ws.Range("I81").Value = 11905810
$ws.Range("J81").Value = 800
$ws.Range("K81").Value = 23811620
$ws.Range("L81").Value = 1600
$ws.Range("M81").Value = -23810559
$ws.Range("N81").Value = -3722
# Row 82
$ws.Range("H82").Value = 43100
$ws.Range("J82").Value = 43100
$ws.Range("L82").Value = 43100
$ws.Range("N82").Value = -43866
# Row 84
$ws.Range("H84").Value = 10205094
$ws.Range("I84").Value = 11905810
$ws.Range("J84").Value = 800
$ws.Range("K84").Value = 119058100
$ws.Range("L84").Value = 8000
$ws.Range("M84").Value = -119052796
$ws.Range("N84").Value = -18608
# Row 85
$ws.Range("H85").Value = 43100
$ws.Range("J85").Value = 43100
$ws.Range("L85").Value = 43100
$ws.Range("N85").Value = -45752
